# Add 2022-Q1 data: turn the old "总计" sheet into the new "2022-Q1"
# fund-holdings detail sheet (reusing its sheetId/position), and append a
# brand-new "总计" sheet after it carrying the refreshed summary table
# (now including the 2022-Q1 row).

function Set-TextCell($ws, $row, $col, $text) {
    # Force text storage (so numeric-looking strings like "011251" or
    # "31.29" keep their original representation instead of being
    # coerced to a Number), then drop back to the workbook's default
    # "Normal" style so no stray NumberFormat/quotePrefix sticks around.
    $ws.Cells.Item($row, $col).Value = "'" + $text
    $ws.Cells.Item($row, $col).Style = "Normal"
}

$wb = $excel.ActiveWorkbook

$q3 = $wb.Worksheets.Item("2021-Q3")
$newSheet = $wb.Worksheets.Item("总计")

# Wipe the old "总计" content/format and repurpose this sheet (keeping its
# sheetId + tab position) as the new "2022-Q1" detail sheet.
$newSheet.Cells.Clear()
$newSheet.Name = "2022-Q1"

# --- formatting -----------------------------------------------------
# Reuse "2021-Q3"'s header-row and first-data-column formatting so the
# new sheet visually matches its siblings.
$q3.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$q3.Range("A2").Copy()
$newSheet.Range("A2:A16").PasteSpecial(-4122)

# --- header row -------------------------------------------------------
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# --- data rows --------------------------------------------------------
$rows = @(
    @("011251", "华安聚嘉精选混合A", "31.29", "89.00", "2.26", "0.7072", 9),
    @("011252", "华安聚嘉精选混合C", "14.89", "89.00", "2.26", "0.3365", 9),
    @("001236", "博时丝路主题股票A", "10.51", "84.58", "2.86", "0.3006", 7),
    @("009740", "博时研究臻选三年持有期灵活配置混合A", "8.15", "81.86", "2.75", "0.2241", 8),
    @("010296", "万家互联互通中国优势量化策略混合A", "5.55", "92.06", "3.78", "0.2098", 3),
    @("011340", "博时战略新材料主题混合A", "1.58", "84.14", "2.84", "0.0449", 9),
    @("290008", "泰信发展主题混合", "0.68", "81.03", "5.03", "0.0342", 7),
    @("010690", "万家互联互通核心资产量化策略混合A", "0.85", "94.05", "4.01", "0.0341", 4),
    @("010297", "万家互联互通中国优势量化策略混合C", "0.53", "92.06", "3.78", "0.0200", 3),
    @("011341", "博时战略新材料主题混合C", "0.63", "84.14", "2.84", "0.0179", 9),
    @("009741", "博时研究臻选三年持有期灵活配置混合C", "0.49", "81.86", "2.75", "0.0135", 8),
    @("004284", "华宝新优选一年定期开放灵活配置混合", "0.64", "38.91", "1.70", "0.0109", 8),
    @("002556", "博时丝路主题股票C", "0.33", "84.58", "2.86", "0.0094", 7),
    @("010691", "万家互联互通核心资产量化策略混合C", "0.20", "94.05", "4.01", "0.0080", 4),
    @("002020", "国都创新驱动灵活配置混合", "0.15", "74.87", "2.67", "0.0040", 3)
)

$r = 2
foreach ($item in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2
    Set-TextCell $newSheet $r 2 $item[0]
    Set-TextCell $newSheet $r 3 $item[1]
    Set-TextCell $newSheet $r 4 $item[2]
    Set-TextCell $newSheet $r 5 $item[3]
    Set-TextCell $newSheet $r 6 $item[4]
    Set-TextCell $newSheet $r 7 $item[5]
    $newSheet.Cells.Item($r, 8).Value = $item[6]
    $r++
}

# --- brand-new "总计" sheet, right after "2022-Q1" ---------------------
$total = $wb.Worksheets.Add($null, $newSheet)
$total.Name = "总计"

# Match the summary sheet's existing look (header + first-data-column
# style) by copying formatting from "2021-Q3" too.
$q3.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$q3.Range("A2").Copy()
$total.Range("A2:A4").PasteSpecial(-4122)

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

$total.Cells.Item(2, 1).Value = 0
Set-TextCell $total 2 2 "2022-Q1"
$total.Cells.Item(2, 3).Value = 15
$total.Cells.Item(2, 4).Value = 1.98

$total.Cells.Item(3, 1).Value = 1
Set-TextCell $total 3 2 "2021-Q3"
$total.Cells.Item(3, 3).Value = 14
$total.Cells.Item(3, 4).Value = 3.85

$total.Cells.Item(4, 1).Value = 2
Set-TextCell $total 4 2 "2021-Q2"
$total.Cells.Item(4, 3).Value = 2
$total.Cells.Item(4, 4).Value = 1.14
